# [#641] Read orientation: call flattentool.unflatten with metatab_vertical_orientation
# (i.e. vertical headings at column 0).
#
# This reshapes the "Meta" tab from a horizontal layout (headers in row 1,
# values in row 2) to a vertical layout (header in column A, value in column
# B) and tidies the extensions URL list (no more stray spaces around the
# `;` separators) plus the accompanying view/formatting tweaks that ship
# with the fixture re-export (90% zoom, wider columns, hidden spacer
# columns on the Meta tab).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: best achievable ColumnWidth input for a desired OOXML width,
# compensating for the engine's pixel-quantised round-trip
# (output = (round(input*6)+5)/6) so the saved width lands as close as
# possible to the real target.
# ---------------------------------------------------------------------
function Set-ColWidth($range, [double]$target) {
    $n = [Math]::Round($target * 6 - 5)
    if ($n -lt 0) { $n = 0 }
    $range.ColumnWidth = ($n / 6)
}

# ---------------------------------------------------------------------
# releases
# ---------------------------------------------------------------------
$wsReleases = $wb.Worksheets.Item("releases")
$wsReleases.Activate()
Set-ColWidth $wsReleases.Columns.Item(1) 17.3238866396761
$excel.ActiveWindow.Zoom = 90
$wsReleases.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# documents
# ---------------------------------------------------------------------
$wsDocuments = $wb.Worksheets.Item("documents")
$wsDocuments.Activate()
Set-ColWidth $wsDocuments.Columns.Item(1) 17.3238866396761
Set-ColWidth $wsDocuments.Columns.Item(2) 41.5101214574899
Set-ColWidth $wsDocuments.Columns.Item(3) 64.5587044534413
Set-ColWidth $wsDocuments.Range("D1").EntireColumn 17.3238866396761
$excel.ActiveWindow.Zoom = 90
$wsDocuments.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# items
# ---------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("items")
$wsItems.Activate()
Set-ColWidth $wsItems.Columns.Item(1) 17.3238866396761
Set-ColWidth $wsItems.Columns.Item(2) 35.7975708502024
Set-ColWidth $wsItems.Columns.Item(3) 39.3279352226721
Set-ColWidth $wsItems.Range("D1").EntireColumn 17.3238866396761
$excel.ActiveWindow.Zoom = 90
$wsItems.Range("B3").Select() | Out-Null

# ---------------------------------------------------------------------
# Meta: transpose version/extensions from a row-pair into column headings,
# clean up the extensions URL, drop the now-redundant hyperlink, resize
# columns (B becomes the wide "value" column, C:D become a hidden spacer)
# and nudge the used range down to row 6 to match the re-exported fixture.
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Meta")
$wsMeta.Activate()

$newExtensionsUrl = "https://raw.githubusercontent.com/open-contracting/ocds_metrics_extension/master/extension.json;https://raw.githubusercontent.com/open-contracting/ocds_extension_parties/master/extension.json;https://raw.githubusercontent.com/open-contracting/ocds_partyDetails_scale_extension/master/extension.json"

# Drop the hyperlink that used to live on B2 (old layout) before we shuffle
# the values around.
$wsMeta.Range("B2").Hyperlinks.Delete()

# Old layout:  A1=version      B1=extensions
#              A2=1.1           B2=<url>
# New layout:  A1=version      B1=1.1
#              A2=extensions   B2=<url, cleaned>
$wsMeta.Range("B1").Value = 1.1
$wsMeta.Range("A2").Value = "extensions"
$wsMeta.Range("B2").Value = $newExtensionsUrl

# B1/B2 lose the Arial "header" style (style index 1) and fall back to the
# plain default font (style index 0).
foreach ($addr in @("B1", "B2")) {
    $cell = $wsMeta.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Color = 0
}

# Column layout: A narrower, B much wider (holds the long value strings),
# C:D collapsed to a hidden spacer, E+ a bit wider too.
Set-ColWidth $wsMeta.Columns.Item(1) 14.3481781376518
Set-ColWidth $wsMeta.Columns.Item(2) 244.708502024291
$wsMeta.Range("C1:D1").EntireColumn.ColumnWidth = 0
$wsMeta.Range("C1:D1").EntireColumn.Hidden = $true
Set-ColWidth $wsMeta.Range("E1").EntireColumn 20.668016194332

# Touch row 6 so the sheet's used range grows to A1:B6, matching the
# re-exported fixture (an otherwise-empty spacer row at the bottom).
$wsMeta.Cells.Item(6, 1).NumberFormat = "General"
$wsMeta.Rows.Item(6).RowHeight = 14.2

$excel.ActiveWindow.Zoom = 90
$wsMeta.Range("B16").Select() | Out-Null
